$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "(318294931, Shalev  Afanasenko: -10,0)"
$ws.Range("B1").Value = "(305487936, Avihai  Kipnis: -7,9)"
$ws.Range("C1").Value = "(313227928, Aviv  Levi: -1,-6)"
$ws.Range("D1").Value = "(205807308, Sariel  Basis: 5,4)"
$ws.Range("E1").Value = "(315891549, Raz  Halaby: 2,9)"
$ws.Range("F1").Value = "(315060103, Dan  Mshelh: 4,8)"
$ws.Range("G1").Value = "(313925141, Elad   Amer: -5,3)"

$ws.Range("A3").Value = "cost: 342.59875744761035"
$ws.Range("A4").Value = "time: 45.37125106394434"
